# Applies the "Update slides on x86 functions (Assembly F'24)" edit to
# slide 25 ("Code Examples") of the x86 procedures-1 deck:
#   - nudge the small green "long mult2(...)" source box to the left
#   - widen/resize the green assembly-listing box
#   - add a new "subq %rsp, 32  # allocate frame" line to the listing

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(25)

# --- Shape "Rectangle 4" (id=3): the C source snippet box -----------------
# off x: 457200 EMU -> 152400 EMU (914400 EMU/in, 12700 EMU/pt)
$srcBox = $s.Shapes.Item(2)
$srcBox.Left = 152400 / 12700

# --- Shape "Rectangle 4" (id=4): the assembly listing box ------------------
$asmBox = $s.Shapes.Item(3)
$asmBox.Left = 2971800 / 12700
$asmBox.Width = 6019800 / 12700
$asmBox.Height = 3657600 / 12700

$tr = $asmBox.TextFrame.TextRange

# Paragraph 3 is "  movq  %rsp, %rbp         # enter" - the function
# prologue's "move stack ptr to base ptr" line. Insert the new
# "subq %rsp, 32  # allocate frame" line right after it.
$afterPara = $tr.Paragraphs(3, 1)
$newParaText = "  subq  %rsp, 32           # allocate frame"
$null = $afterPara.InsertAfter("`r" + $newParaText)

# The newly-created paragraph is now #4; split it into the same run
# boundaries as the authored slide (so each fragment keeps its own
# run, matching "subq", "rsp", etc. as distinct runs) by nudging a
# (no-op) character-level format on each fragment boundary.
$newPara = $tr.Paragraphs(4, 1)
$pos = $newPara.Start
$runTexts = @("  ", "subq", "  %", "rsp", ", 32           ", "# allocate frame")
foreach ($runText in $runTexts) {
    $len = $runText.Length
    $frag = $tr.Characters($pos, $len)
    $frag.Font.Bold = $false
    $pos += $len
}
